# Add "2022-Q3" sheet data.
#
# The workbook has a "总计" (totals) sheet followed by one sheet per
# quarter (most-recent-first). We insert a brand new "2022-Q3" sheet
# right after "总计" (pushing every existing quarter sheet down by one
# tab, unchanged), populate it with the new quarter's fund data, and add
# a corresponding summary row on "总计".

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q3" worksheet ------------------------------
# Duplicate the current "2022-Q2" sheet (position 2) right after "总计"
# (position 1) so the new sheet inherits identical formatting/styles,
# then rename it and overwrite its data with the Q3 numbers.
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($null, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Fund codes (column B) and row indices (column A) stay as copied from
# 2022-Q2 (009562 / 009563 / 486001); only names + figures change.
$q3Sheet.Cells.Item(2, 3).Value = "工银全球股票（QDII）美元"
$q3Sheet.Cells.Item(2, 4).Value = "5.89"
$q3Sheet.Cells.Item(2, 5).Value = "93.72"
$q3Sheet.Cells.Item(2, 6).Value = "1.90"
$q3Sheet.Cells.Item(2, 7).Value = "0.1119"
$q3Sheet.Cells.Item(2, 8).Value = 9

$q3Sheet.Cells.Item(3, 3).Value = "工银全球股票（QDII）港币"
$q3Sheet.Cells.Item(3, 4).Value = "5.89"
$q3Sheet.Cells.Item(3, 5).Value = "93.72"
$q3Sheet.Cells.Item(3, 6).Value = "1.90"
$q3Sheet.Cells.Item(3, 7).Value = "0.1119"
$q3Sheet.Cells.Item(3, 8).Value = 9

$q3Sheet.Cells.Item(4, 3).Value = "工银瑞信中国机会全球配置股票（QDII）人民币"
$q3Sheet.Cells.Item(4, 4).Value = "5.89"
$q3Sheet.Cells.Item(4, 5).Value = "93.72"
$q3Sheet.Cells.Item(4, 6).Value = "1.90"
$q3Sheet.Cells.Item(4, 7).Value = "0.1119"
$q3Sheet.Cells.Item(4, 8).Value = 9

# --- 2. Update the "总计" summary sheet ----------------------------------
# Insert a new row for 2022-Q3 right under the header, re-numbering the
# trailing index column (A) and shifting every other quarter down a row.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.34

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(3, 3).Value = 3
$totalSheet.Cells.Item(3, 4).Value = 0.49

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(4, 3).Value = 3
$totalSheet.Cells.Item(4, 4).Value = 0.4

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(5, 3).Value = 3
$totalSheet.Cells.Item(5, 4).Value = 0.78

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(6, 3).Value = 3
$totalSheet.Cells.Item(6, 4).Value = 0.3

$totalSheet.Cells.Item(7, 1).Value = 5
$totalSheet.Cells.Item(7, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(7, 3).Value = 1
$totalSheet.Cells.Item(7, 4).Value = 0
